$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new credential row (fill order matches the original authoring
# order so new shared strings are appended in the same sequence).
$ws.Range("C5").Value = "profileSettings"
$ws.Range("B5").Value = "automation+1@thinkbridge.in"
$ws.Range("A5").Value = "Neha Automation "
$ws.Range("D5").Value = "Consero234$"

# Turn the email cell into a hyperlink, then restore the "Hyperlink" look
# (matching the styling used by the other rows' email cells).
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:automation+1@thinkbridge.in") | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null

# Leave the same selection state recorded by the author.
$ws.Range("E5").Select() | Out-Null
